# Update test results data - DU FBS Mock 2
# Generated with VH Results Processing System

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 22 corresponds to student ID 7842103 (Shafqat Nur).
# Update the raw MCQ correct/wrong inputs; all dependent formula
# cells (marks, percentages, totals, and RANK.EQ-based ranks across
# the whole sheet) recalculate automatically.
$ws.Range("C22").Value = 15   # English Correct
$ws.Range("D22").Value = 1    # English Wrong
$ws.Range("G22").Value = 10   # Adv English Correct
$ws.Range("H22").Value = 2    # Adv English Wrong
$ws.Range("K22").Value = 8    # Business Studies Correct
$ws.Range("L22").Value = 8    # Business Studies Wrong
$ws.Range("S22").Value = 9    # Economics Correct
$ws.Range("T22").Value = 6    # Economics Wrong

# Move the active cell selection to Y22, matching the latest editing
# location in the source file.
$ws.Range("Y22").Select()
